$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Workbook-level: switch reference style back to A1 (drop refMode="R1C1")
$excel.ReferenceStyle = 1

# ---- Sheet1 "TO BE FILLED": add new column AA "PreferredDeliveryDate" ----
$ws1.Range("AA1").Value = "PreferredDeliveryDate"
$ws1.Range("D1").Copy()
$ws1.Range("AA1").PasteSpecial(-4122)
$ws1.Columns.Item(27).ColumnWidth = 23.6667

# ---- Sheet2 "FILLING RULES AND EXAMPLE": add new column AA ----
$ws2.Range("AA1").Value = "PreferredDeliveryDate"
$ws2.Range("D1").Copy()
$ws2.Range("AA1").PasteSpecial(-4122)
$ws2.Columns.Item(27).ColumnWidth = 18.6667

$ws2.Range("AA2").Value = "date when the order should be delivered"
$ws2.Range("D2").Copy()
$ws2.Range("AA2").PasteSpecial(-4122)

$ws2.Range("AA3").Value = "not required"
$ws2.Range("D3").Copy()
$ws2.Range("AA3").PasteSpecial(-4122)

# Example delivery date for the second sample order (row 5)
$ws2.Range("AA5").Value = 45533
$ws2.Range("N5").Copy()
$ws2.Range("AA5").PasteSpecial(-4122)
$ws2.Range("AA5").NumberFormat = "mm-dd-yy"

$excel.CutCopyMode = $false

# ---- Selection / view bookkeeping ----
$ws2.Activate()
$ws2.Range("AA6").Select()
$ws1.Activate()
$ws1.Range("A2").Select()
